# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets contain identical data, and both need the same updates.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new F-column value
$updates = @{
    2  = 1948
    4  = 118
    7  = 1622
    9  = 641
    16 = 136
    17 = 110
    19 = 3761
    21 = 15
    23 = 346
    25 = 415
    28 = 1547
    29 = 12
    30 = 149
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
